$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '47.166.73'
$ws.Range("E2").Value = '  +1.19%  '

# Row 3
$ws.Range("D3").Value = '2.488.71'

# Row 4
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.48%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.08%  '

# Row 7
$ws.Range("E7").Value = '  -0.06%  '

# Row 8
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.533'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.10%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.77'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.10%  '

# Row 11
$ws.Range("E11").Value = '  -0.82%  '

# Row 12
$ws.Range("E12").Value = '  +0.25%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.24'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.46%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.12'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.41%  '

# Row 15
$ws.Range("D15").Value = '2.877.82'
$ws.Range("E15").Value = '  +0.83%  '

# Row 16
$ws.Range("D16").Value = '2.494.27'
$ws.Range("E16").Value = '  -0.79%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.844'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.16%  '

# Row 18
$ws.Range("D18").Value = '47.075.96'
$ws.Range("E18").Value = '  +1.25%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.72'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.09%  '

# Row 20
$ws.Range("E20").Value = '  +2.09%  '

# Row 21
$ws.Range("B21").Value = 'ImmutableX'
$ws.Range("C21").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +16.29%  '

# Row 22
$ws.Range("B22").Value = 'ShibaInu'
$ws.Range("C22").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D22").Value = '0.0₃0934'
$ws.Range("E22").Value = '  -0.42%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.64'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.18%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '245.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.14%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.56'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.55%  '

# Row 26
$ws.Range("E26").Value = '  -0.02%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.65%  '

# Row 28
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.04'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.34%  '

# Row 29
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.140'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +8.42%  '

# Row 30
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.06'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.77%  '

# Row 31
$ws.Range("B31").Value = 'OKB'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '49.91'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.61%  '

# Row 32
$ws.Range("B32").Value = 'Celestia'
$ws.Range("C32").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.98'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.85%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.35'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.36%  '

# Row 34
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0783'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.90%  '

# Row 35
$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.01'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.30%  '

# Row 36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.95'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.14%  '

# Row 37
$ws.Range("E37").Value = '  +0.17%  '

# Row 38
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.94'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.05%  '

# Row 39
$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.112'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.16%  '

# Row 40
$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.22'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.05%  '

# Row 41
$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '119.02'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.90%  '

# Row 42
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.30'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.06%  '

# Row 43
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0294'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.14%  '

# Row 44
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '1.979.01'
$ws.Range("E44").Value = '  -0.37%  '

# Row 45
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.02'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.08%  '

# Row 46
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.01'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.67%  '

# Row 47
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.05'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.23%  '

# Row 48
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.77'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.84%  '

# Row 49
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.12'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.29%  '

# Row 50
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '57.09'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.23%  '

# Row 51
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '76.92'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.80%  '

Write-Output "Updated cryptos list"